$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the bold/bordered
# header style already used by the other header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF) for rows 2-21.
$iValues = @(9, 8, 5, 5, 6, 9, 8, 7, 7, 8, 7, 7, 8, 9, 5, 6, 8, 7, 7, 7)
$jValues = @(9, 8, 5, 6, 7, 9, 8, 7, 8, 8, 7, 7, 8, 9, 6, 8, 8, 8, 7, 7)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
